$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '20.557.23'
$ws.Range("E2").Value = '  +1.84%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.472.93'
$ws.Range("E3").Value = '  +2.81%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  -0.42%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9619'
$ws.Range("E5").Value = '  +5.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '277.48'
$ws.Range("E6").Value = '  +0.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3625'
$ws.Range("E7").Value = '  -0.37%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3077'
$ws.Range("E8").Value = '  -0.45%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '39.72'
$ws.Range("E9").Value = '  +2.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.073'
$ws.Range("E10").Value = '  +5.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06648'
$ws.Range("E11").Value = '  +2.33%  '
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.512'
$ws.Range("E13").Value = '  +3.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.24'
$ws.Range("E14").Value = '  +4.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.9624'
$ws.Range("E15").Value = '  +2.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.164'
$ws.Range("E16").Value = '  +2.37%  '
$ws.Range("E17").Value = '  +1.59%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.473.88'
$ws.Range("E18").Value = '  +2.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.05946'
$ws.Range("E19").Value = '  +5.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.82'
$ws.Range("E20").Value = '  +2.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.504'
$ws.Range("E21").Value = '  +2.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.54'
$ws.Range("E22").Value = '  +1.83%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.18'
$ws.Range("E23").Value = '  +4.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.254'
$ws.Range("E24").Value = '  +1.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '20.567.74'
$ws.Range("E25").Value = '  +1.55%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.12'
$ws.Range("E26").Value = '  +3.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.131'
$ws.Range("E27").Value = '  -0.44%  '
$ws.Range("E28").Value = '  +2.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.636.20'
$ws.Range("E29").Value = '  +2.78%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '113.87'
$ws.Range("E30").Value = '  +4.35%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.886'
$ws.Range("E31").Value = '  +0.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.944'
$ws.Range("E32").Value = '  +4.34%  '
$ws.Range("E33").Value = '  +4.73%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.8053'
$ws.Range("E34").Value = '  +1.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.509'
$ws.Range("E35").Value = '  +4.90%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.213'
$ws.Range("E36").Value = '  +7.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05767'
$ws.Range("E37").Value = '  -2.45%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.730'
$ws.Range("E38").Value = '  +2.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02051'
$ws.Range("E39").Value = '  +3.80%  '
$ws.Range("E40").Value = '  +4.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.41'
$ws.Range("E41").Value = '  +3.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1878'
$ws.Range("E42").Value = '  +2.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.420'
$ws.Range("E43").Value = '  +5.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5281'
$ws.Range("E44").Value = '  +1.66%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.24'
$ws.Range("E45").Value = '  +2.20%  '
$ws.Range("E46").Value = '  +0.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '118.87'
$ws.Range("E47").Value = '  +0.70%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5203'
$ws.Range("E48").Value = '  +2.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.817'
$ws.Range("E49").Value = '  +4.21%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06464'
$ws.Range("E50").Value = '  +2.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9879'
$ws.Range("E51").Value = '  -0.13%  '
